$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as plain text (matching the existing rows),
# not Excel's autodetected date serials. Force text storage by formatting
# the cell as Text before assigning the value, then clear the formatting
# override so the cell's style stays at the sheet default (same as the
# other rows in the table).
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "12/30/2025"
$ws.Range("A36").ClearFormats()

$ws.Range("B36").Value = 11954
$ws.Range("C36").Value = 0.2155915824299713
$ws.Range("D36").Value = 0.7844084175700287
$ws.Range("E36").Value = -148.32
$ws.Range("F36").Value = -27.72
$ws.Range("G36").Value = -21214.94
$ws.Range("H36").Value = -69.34999999999999
$ws.Range("I36").Value = -475.67
$ws.Range("J36").Value = -15.58
